$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bksr")

# Convert the B column (年份/time) from text labels ("2023年","2024年","2025年")
# to plain numeric years (2023, 2024, 2025)
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 2).Value = 2023
}
for ($r = 10; $r -le 17; $r++) {
    $ws.Cells.Item($r, 2).Value = 2024
}
for ($r = 18; $r -le 26; $r++) {
    $ws.Cells.Item($r, 2).Value = 2025
}

# Move the active selection
$ws.Range("I23").Select()
